$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells stay as Text, matching the source data which
# stores numeric-looking prices (e.g. "93.07", "71.20") as strings, not numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '41.693.24'
$ws.Range("E2").Value = '  +0.41%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.475.91'
$ws.Range("E3").Value = '  +0.39%  '

$ws.Range("E4").Value = '  +0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.71'
$ws.Range("E5").Value = '  +1.29%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '93.07'
$ws.Range("E6").Value = '  +1.74%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.553'
$ws.Range("E7").Value = '  +1.14%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("E9").Value = '  +1.06%  '

$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0862'
$ws.Range("E10").Value = '  +9.52%  '

$ws.Range("B11").Value = 'Avalanche'
$ws.Range("C11").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '33.13'
$ws.Range("E11").Value = '  +3.17%  '

$ws.Range("E12").Value = '  +0.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.857.84'
$ws.Range("E13").Value = '  +0.39%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.89'
$ws.Range("E14").Value = '  +0.84%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.80'
$ws.Range("E15").Value = '  -1.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.499.02'
$ws.Range("E16").Value = '  +0.84%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.788'
$ws.Range("E17").Value = '  +3.07%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '41.687.91'
$ws.Range("E18").Value = '  +0.47%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.48'
$ws.Range("E19").Value = '  +0.02%  '

$ws.Range("E20").Value = '  +0.55%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '71.20'
$ws.Range("E21").Value = '  -0.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.32'
$ws.Range("E22").Value = '  +2.43%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '239.77'
$ws.Range("E23").Value = '  +1.94%  '

$ws.Range("E24").Value = '  +1.31%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.93'
$ws.Range("E25").Value = '  +2.48%  '

$ws.Range("E26").Value = '  +0.01%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '24.76'
$ws.Range("E27").Value = '  +0.98%  '

$ws.Range("E28").Value = '  +2.49%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.82'
$ws.Range("E29").Value = '  +1.88%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.03'
$ws.Range("E30").Value = '  +2.20%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '158.34'
$ws.Range("E31").Value = '  +1.66%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.53'
$ws.Range("E32").Value = '  +1.87%  '

$ws.Range("E33").Value = '  -0.12%  '

$ws.Range("E34").Value = '  +0.86%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0769'
$ws.Range("E35").Value = '  +1.94%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.32'
$ws.Range("E36").Value = '  +1.60%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.88'
$ws.Range("E37").Value = '  +5.15%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.93'
$ws.Range("E38").Value = '  +2.48%  '

$ws.Range("E39").Value = '  +2.05%  '

$ws.Range("E40").Value = '  +0.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.03'
$ws.Range("E41").Value = '  +1.08%  '

$ws.Range("E42").Value = '  +7.37%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.994.37'
$ws.Range("E43").Value = '  +2.86%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '19.37'
$ws.Range("E44").Value = '  +4.89%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0286'
$ws.Range("E45").Value = '  +1.19%  '

$ws.Range("E46").Value = '  +2.61%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.33'
$ws.Range("E47").Value = '  +3.35%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.714.23'
$ws.Range("E48").Value = '  +0.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '97.44'
$ws.Range("E49").Value = '  +0.57%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '74.34'
$ws.Range("E50").Value = '  +4.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '67.20'
$ws.Range("E51").Value = '  +0.72%  '
